# kapital_bank_excel 2020_Q3 capital_adequacy: add the "Table_2" worksheet
# (norms / actual ratios table) and tidy up a few stray empty cells left on
# Table_1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Table_1: drop the empty placeholder cells (B2, A3, B37) -------------
$ws1.Range("B2").ClearContents()
$ws1.Range("A3").ClearContents()
$ws1.Range("B37").ClearContents()

# --- add the new "Table_2" worksheet right after Table_1 -----------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Table_2"

# Match Table_1's sheet-level look & feel (outline summary position + page
# margins) on the freshly added sheet.
$ws2.Outline.SummaryColumn = 1
$ws2.Outline.SummaryRow = 1
$ws2.PageSetup.LeftMargin   = $ws1.PageSetup.LeftMargin
$ws2.PageSetup.RightMargin  = $ws1.PageSetup.RightMargin
$ws2.PageSetup.TopMargin    = $ws1.PageSetup.TopMargin
$ws2.PageSetup.BottomMargin = $ws1.PageSetup.BottomMargin
$ws2.PageSetup.HeaderMargin = $ws1.PageSetup.HeaderMargin
$ws2.PageSetup.FooterMargin = $ws1.PageSetup.FooterMargin

# Match the bold/bordered/centered header style used on Table_1's header row.
$ws1.Range("A1:B1").Copy()
$ws2.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- header row ------------------------------------------------------------
$ws2.Range("A1").Value = "Əmsal"
$ws2.Range("B1").Value = "Norma (Sistem əhəmiyyətli)"
$ws2.Range("C1").Value = "Norma (Banklar istisna)"
$ws2.Range("D1").Value = "Fakt"

# --- data rows ---------------------------------------------------------
# Percentage-looking text ("6.0%", "5.19%", ...) has to stay literal text
# (not get auto-converted into a numeric percentage by Excel's input
# parser), so those cells are pre-formatted as Text before the value is
# written.
$percentCells = @("B2","C2","D2","B3","C3","D3","D4")
foreach ($addr in $percentCells) {
    $ws2.Range($addr).NumberFormat = "@"
}

$ws2.Range("A2").Value = "9.  I dərəcəli  kapitalın  adekvatlıq əmsalı"
$ws2.Range("B2").Value = "6.0%"
$ws2.Range("C2").Value = "5.0%"
$ws2.Range("D2").Value = "10.49%"

$ws2.Range("A3").Value = "10. məcmu kapitalın  adekvatlıq  əmsalı"
$ws2.Range("B3").Value = "11.0%"
$ws2.Range("C3").Value = "9.0%"
$ws2.Range("D3").Value = "20.6%"

$ws2.Range("A4").Value = "11. Leverec əmsalı"
$ws2.Range("B4").Value = "minimum 5%"
$ws2.Range("C4").Value = "minimum 4%"
$ws2.Range("D4").Value = "5.19%"

$ws1.Select() | Out-Null
$ws1.Range("A1").Select() | Out-Null
